$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 917100
$ws.Range("L2").Value = 752900
$ws.Range("M2").Value = 2024
$ws.Range("N2").Value = "None"
$ws.Range("O2").Value = 2148
$ws.Range("P2").Value = 1964
$ws.Range("Q2").Value = 0.4699954086317722
$ws.Range("R2").Value = "Acres"
$ws.Range("S2").Value = "Forced air,Gas"
$ws.Range("T2").Value = "None"

# Row 3
$ws.Range("K3").Value = 1216400
$ws.Range("L3").Value = 931700
$ws.Range("M3").Value = 2024
$ws.Range("N3").Value = "'2012-07-24"
$ws.Range("O3").Value = 2258
$ws.Range("P3").Value = 1929
$ws.Range("Q3").Value = 7590
$ws.Range("R3").Value = "Square Feet"
$ws.Range("S3").Value = "Electric Baseboard,Steam,Natural Gas"
$ws.Range("T3").Value = "None"

# Row 4
$ws.Range("K4").Value = 1035600
$ws.Range("L4").Value = 831600
$ws.Range("M4").Value = 2024
$ws.Range("N4").Value = "'2008-06-27"
$ws.Range("O4").Value = 1668
$ws.Range("P4").Value = 1952
$ws.Range("Q4").Value = 0.5799816345270891
$ws.Range("R4").Value = "Acres"
$ws.Range("S4").Value = "Forced air,Gas"
$ws.Range("T4").Value = "None"

# Row 5
$ws.Range("K5").Value = "None"
$ws.Range("L5").Value = 995700
$ws.Range("M5").Value = 2024
$ws.Range("N5").Value = "'1999-04-09"
$ws.Range("O5").Value = "None"
$ws.Range("P5").Value = "None"
$ws.Range("Q5").Value = 0.9199954086317723
$ws.Range("R5").Value = "Acres"
$ws.Range("S5").Value = "None"
$ws.Range("T5").Value = "None"

# Row 6
$ws.Range("K6").Value = 652800
$ws.Range("L6").Value = 505800
$ws.Range("M6").Value = 2024
$ws.Range("N6").Value = "'1987-09-01"
$ws.Range("O6").Value = 2366
$ws.Range("P6").Value = 1987
$ws.Range("Q6").Value = 0.55
$ws.Range("R6").Value = "Acres"
$ws.Range("S6").Value = "Other,Oil"
$ws.Range("T6").Value = "None"

# Row 7
$ws.Range("K7").Value = 761800
$ws.Range("L7").Value = 622300
$ws.Range("M7").Value = 2024
$ws.Range("N7").Value = "'2008-10-16"
$ws.Range("O7").Value = 1612
$ws.Range("P7").Value = 2000
$ws.Range("Q7").Value = 9583
$ws.Range("R7").Value = "Square Feet"
$ws.Range("S7").Value = "Forced air,Gas"
$ws.Range("T7").Value = "Central"

# Row 8
$ws.Range("K8").Value = 805200
$ws.Range("L8").Value = 655300
$ws.Range("M8").Value = 2024
$ws.Range("N8").Value = "'1996-05-24"
$ws.Range("O8").Value = 1333
$ws.Range("P8").Value = 1930
$ws.Range("Q8").Value = 0.3994490358126722
$ws.Range("R8").Value = "Acres"
$ws.Range("S8").Value = "Heat pump,Electric,Solar"
$ws.Range("T8").Value = "Central,Solar"

# Row 9
$ws.Range("K9").Value = 907000
$ws.Range("L9").Value = 736800
$ws.Range("M9").Value = 2024
$ws.Range("N9").Value = "None"
$ws.Range("O9").Value = 2520
$ws.Range("P9").Value = 1989
$ws.Range("Q9").Value = 1.069995408631772
$ws.Range("R9").Value = "Acres"
$ws.Range("S9").Value = "Other,Oil"
$ws.Range("T9").Value = "None"

# Row 10
$ws.Range("K10").Value = 1219200
$ws.Range("L10").Value = 985500
$ws.Range("M10").Value = 2024
$ws.Range("N10").Value = "'2015-08-14"
$ws.Range("O10").Value = 3564
$ws.Range("P10").Value = 2000
$ws.Range("Q10").Value = 4.479981634527089
$ws.Range("R10").Value = "Acres"
$ws.Range("S10").Value = "Forced air,Heat pump,Stove,Oil,Solar,Wood / Pellet"
$ws.Range("T10").Value = "Central,Solar"

# Row 11
$ws.Range("K11").Value = 409300
$ws.Range("L11").Value = 332300
$ws.Range("M11").Value = 2024
$ws.Range("N11").Value = "None"
$ws.Range("O11").Value = 1958
$ws.Range("P11").Value = 1950
$ws.Range("Q11").Value = 0.2899908172635445
$ws.Range("R11").Value = "Acres"
$ws.Range("S11").Value = "Forced air,Oil"
$ws.Range("T11").Value = "None"

# Row 12
$ws.Range("K12").Value = 486900
$ws.Range("L12").Value = 422700
$ws.Range("M12").Value = 2024
$ws.Range("N12").Value = "'2017-08-21"
$ws.Range("O12").Value = 1656
$ws.Range("P12").Value = 1994
$ws.Range("Q12").Value = 0.6699954086317723
$ws.Range("R12").Value = "Acres"
$ws.Range("S12").Value = "Forced air,Oil"
$ws.Range("T12").Value = "None"

# Row 13
$ws.Range("K13").Value = 698200
$ws.Range("L13").Value = 557400
$ws.Range("M13").Value = 2024
$ws.Range("N13").Value = "'2013-09-16"
$ws.Range("O13").Value = 1356
$ws.Range("P13").Value = 1954
$ws.Range("Q13").Value = 10018
$ws.Range("R13").Value = "Square Feet"
$ws.Range("S13").Value = "Gas"
$ws.Range("T13").Value = "None"

# Row 14
$ws.Range("K14").Value = 792800
$ws.Range("L14").Value = 680500
$ws.Range("M14").Value = 2024
$ws.Range("N14").Value = "'2017-03-10"
$ws.Range("O14").Value = 3265
$ws.Range("P14").Value = 1999
$ws.Range("Q14").Value = 1.089990817263545
$ws.Range("R14").Value = "Acres"
$ws.Range("S14").Value = "Forced air,Stove,Gas"
$ws.Range("T14").Value = "Central"

# Row 15
$ws.Range("K15").Value = 4478700
$ws.Range("L15").Value = 3764200
$ws.Range("M15").Value = 2024
$ws.Range("N15").Value = "'2016-10-28"
$ws.Range("O15").Value = 2672
$ws.Range("P15").Value = 1996
$ws.Range("Q15").Value = 3.569995408631772
$ws.Range("R15").Value = "Acres"
$ws.Range("S15").Value = "Forced air,Gas"
$ws.Range("T15").Value = "None"

# Row 16
$ws.Range("K16").Value = 1559500
$ws.Range("L16").Value = 1099000
$ws.Range("M16").Value = 2024
$ws.Range("N16").Value = "'2001-08-06"
$ws.Range("O16").Value = 2612
$ws.Range("P16").Value = 1974
$ws.Range("Q16").Value = 0.2599862258953168
$ws.Range("R16").Value = "Acres"
$ws.Range("S16").Value = "Forced air,Gas"
$ws.Range("T16").Value = "None"

# Row 17
$ws.Range("K17").Value = 1916600
$ws.Range("L17").Value = 1622300
$ws.Range("M17").Value = 2024
$ws.Range("N17").Value = "'2006-07-14"
$ws.Range("O17").Value = 4964
$ws.Range("P17").Value = 1991
$ws.Range("Q17").Value = 0.5
$ws.Range("R17").Value = "Acres"
$ws.Range("S17").Value = "Electric"
$ws.Range("T17").Value = "Central"

# Row 18
$ws.Range("K18").Value = 544200
$ws.Range("L18").Value = 396100
$ws.Range("M18").Value = 2024
$ws.Range("N18").Value = "'2008-03-19"
$ws.Range("O18").Value = 1040
$ws.Range("P18").Value = 1958
$ws.Range("Q18").Value = 10018
$ws.Range("R18").Value = "Square Feet"
$ws.Range("S18").Value = "Other,Oil"
$ws.Range("T18").Value = "None"

# Row 19
$ws.Range("K19").Value = 782600
$ws.Range("L19").Value = 645600
$ws.Range("M19").Value = 2024
$ws.Range("N19").Value = "None"
$ws.Range("O19").Value = 2936
$ws.Range("P19").Value = 1992
$ws.Range("Q19").Value = 1.839990817263545
$ws.Range("R19").Value = "Acres"
$ws.Range("S19").Value = "Oil"
$ws.Range("T19").Value = "None"

# Row 20
$ws.Range("K20").Value = 1358300
$ws.Range("L20").Value = 1153100
$ws.Range("M20").Value = 2024
$ws.Range("N20").Value = "None"
$ws.Range("O20").Value = 2020
$ws.Range("P20").Value = 1953
$ws.Range("Q20").Value = 1.109986225895317
$ws.Range("R20").Value = "Acres"
$ws.Range("S20").Value = "Oil"
$ws.Range("T20").Value = "None"

# Row 21
$ws.Range("K21").Value = 3081500
$ws.Range("L21").Value = 1727200
$ws.Range("M21").Value = 2024
$ws.Range("N21").Value = "'2018-09-20"
$ws.Range("O21").Value = 2844
$ws.Range("P21").Value = 2000
$ws.Range("Q21").Value = 0.463682277318641
$ws.Range("R21").Value = "Acres"
$ws.Range("S21").Value = "Solar"
$ws.Range("T21").Value = "None"

# Row 22
$ws.Range("K22").Value = 497900
$ws.Range("L22").Value = 429200
$ws.Range("M22").Value = 2024
$ws.Range("N22").Value = "None"
$ws.Range("O22").Value = 1714
$ws.Range("P22").Value = 1961
$ws.Range("Q22").Value = 1.5
$ws.Range("R22").Value = "Acres"
$ws.Range("S22").Value = "Other,Electric"
$ws.Range("T22").Value = "None"

# Row 23
$ws.Range("K23").Value = 1900300
$ws.Range("L23").Value = 1426800
$ws.Range("M23").Value = 2024
$ws.Range("N23").Value = "'2006-09-20"
$ws.Range("O23").Value = 2849
$ws.Range("P23").Value = 2007
$ws.Range("Q23").Value = 0.4683195592286502
$ws.Range("R23").Value = "Acres"
$ws.Range("S23").Value = "Baseboard,Gas"
$ws.Range("T23").Value = "Central"

# Row 24
$ws.Range("K24").Value = 670200
$ws.Range("L24").Value = 534500
$ws.Range("M24").Value = 2024
$ws.Range("N24").Value = "'2016-03-29"
$ws.Range("O24").Value = 1858
$ws.Range("P24").Value = 1968
$ws.Range("Q24").Value = 0.6899908172635445
$ws.Range("R24").Value = "Acres"
$ws.Range("S24").Value = "Gas"
$ws.Range("T24").Value = "None"

# Row 25
$ws.Range("K25").Value = 1235700
$ws.Range("L25").Value = 922500
$ws.Range("M25").Value = 2024
$ws.Range("N25").Value = "'2009-08-19"
$ws.Range("O25").Value = 3541
$ws.Range("P25").Value = 1971
$ws.Range("Q25").Value = 0.9366620752984389
$ws.Range("R25").Value = "Acres"
$ws.Range("S25").Value = "Other,Gas"
$ws.Range("T25").Value = "None"

# Row 26
$ws.Range("K26").Value = 621400
$ws.Range("L26").Value = 462900
$ws.Range("M26").Value = 2024
$ws.Range("N26").Value = "'2005-06-14"
$ws.Range("O26").Value = 2342
$ws.Range("P26").Value = 1987
$ws.Range("Q26").Value = 0.2514003673094582
$ws.Range("R26").Value = "Acres"
$ws.Range("S26").Value = "Other,Gas"
$ws.Range("T26").Value = "Central,Solar"

# Row 210
$ws.Range("K210").Value = 536900
$ws.Range("L210").Value = 400700
$ws.Range("M210").Value = 2024
$ws.Range("N210").Value = "'2022-06-17"
$ws.Range("O210").Value = 1716
$ws.Range("P210").Value = 1930
$ws.Range("Q210").Value = 6969.6
$ws.Range("R210").Value = "Square Feet"
$ws.Range("S210").Value = "Electric Baseboard,Steam,Natural Gas"
$ws.Range("T210").Value = "None"

Write-Output "Applied K2:T26 and K210:T210 updates"
